$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header label (unnamed: 1_level_1 -> total)
$ws.Range("B2").Value = "total"

# Remove the two trailing rows (old goias/distrito federal slots no longer present)
$ws.Rows("39:40").Delete()

# Rewrite labels (col A) and numeric values (cols B:I) for rows 4-38
$ws.Range("A4").Value = "               brasil"
$ws.Range("B4").Value = 0.46
$ws.Range("C4").Value = 1.63
$ws.Range("D4").Value = 0.75
$ws.Range("E4").Value = 1.07
$ws.Range("F4").Value = 1.04
$ws.Range("G4").Value = 0.87
$ws.Range("H4").Value = 1.76
$ws.Range("I4").Value = 1.92

$ws.Range("A5").Value = "urbana"
$ws.Range("B5").Value = 0.72
$ws.Range("C5").Value = 1.73
$ws.Range("D5").Value = 0.95
$ws.Range("E5").Value = 1.23
$ws.Range("F5").Value = 1.23
$ws.Range("G5").Value = 1.01
$ws.Range("H5").Value = 1.86
$ws.Range("I5").Value = 2

$ws.Range("A6").Value = "rural"
$ws.Range("B6").Value = 2.88
$ws.Range("C6").Value = 3.89
$ws.Range("D6").Value = 2.93
$ws.Range("E6").Value = 4.09
$ws.Range("F6").Value = 4.06
$ws.Range("G6").Value = 4.05
$ws.Range("H6").Value = 7.34
$ws.Range("I6").Value = 7.21

$ws.Range("A7").Value = "norte"
$ws.Range("B7").Value = 2.71
$ws.Range("C7").Value = 5.39
$ws.Range("D7").Value = 3.2
$ws.Range("E7").Value = 4.21
$ws.Range("F7").Value = 4.92
$ws.Range("G7").Value = 2.97
$ws.Range("H7").Value = 4.79
$ws.Range("I7").Value = 5.05

$ws.Range("A8").Value = "rondônia"
$ws.Range("B8").Value = 2.6
$ws.Range("C8").Value = 7.36
$ws.Range("D8").Value = 3.91
$ws.Range("E8").Value = 7.89
$ws.Range("F8").Value = 7.57
$ws.Range("G8").Value = 4.86
$ws.Range("H8").Value = 12.09
$ws.Range("I8").Value = 11.43

$ws.Range("A9").Value = "acre"
$ws.Range("B9").Value = 4.68
$ws.Range("C9").Value = 8.16
$ws.Range("D9").Value = 6.87
$ws.Range("E9").Value = 7.65
$ws.Range("F9").Value = 11.57
$ws.Range("G9").Value = 8.619999999999999
$ws.Range("H9").Value = 13.82
$ws.Range("I9").Value = 15.86

$ws.Range("A10").Value = "amazonas"
$ws.Range("B10").Value = 3.27
$ws.Range("C10").Value = 14
$ws.Range("D10").Value = 3.39
$ws.Range("E10").Value = 6.74
$ws.Range("F10").Value = 7.4
$ws.Range("G10").Value = 4.86
$ws.Range("H10").Value = 10.27
$ws.Range("I10").Value = 11.23

$ws.Range("A11").Value = "roraima"
$ws.Range("B11").Value = 3.45
$ws.Range("C11").Value = 9.81
$ws.Range("D11").Value = 5.54
$ws.Range("E11").Value = 11.98
$ws.Range("F11").Value = 12.44
$ws.Range("G11").Value = 8.82
$ws.Range("H11").Value = 16.36
$ws.Range("I11").Value = 18.93

$ws.Range("A12").Value = "pará"
$ws.Range("B12").Value = 6.41
$ws.Range("C12").Value = 9.93
$ws.Range("D12").Value = 6.77
$ws.Range("E12").Value = 9.140000000000001
$ws.Range("F12").Value = 11.1
$ws.Range("G12").Value = 7.09
$ws.Range("H12").Value = 9.720000000000001
$ws.Range("I12").Value = 11.54

$ws.Range("A13").Value = "amapá"
$ws.Range("B13").Value = 7.42
$ws.Range("C13").Value = 13.11
$ws.Range("D13").Value = 14.43
$ws.Range("E13").Value = 8.300000000000001
$ws.Range("F13").Value = 14.33
$ws.Range("G13").Value = 9.210000000000001
$ws.Range("H13").Value = 17.61
$ws.Range("I13").Value = 17.61

$ws.Range("A14").Value = "tocantins"
$ws.Range("B14").Value = 3.47
$ws.Range("C14").Value = 9.640000000000001
$ws.Range("D14").Value = 5.3
$ws.Range("E14").Value = 7.59
$ws.Range("F14").Value = 7.08
$ws.Range("G14").Value = 5.07
$ws.Range("H14").Value = 11.84
$ws.Range("I14").Value = 9.44

$ws.Range("A15").Value = "nordeste"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 2.74
$ws.Range("D15").Value = 1.58
$ws.Range("E15").Value = 2.44
$ws.Range("F15").Value = 2.04
$ws.Range("G15").Value = 2.36
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 4.29

$ws.Range("A16").Value = "maranhão"
$ws.Range("B16").Value = 3.52
$ws.Range("C16").Value = 10.64
$ws.Range("D16").Value = 7.94
$ws.Range("E16").Value = 12.52
$ws.Range("F16").Value = 8.02
$ws.Range("G16").Value = 12.7
$ws.Range("H16").Value = 16.51
$ws.Range("I16").Value = 21.64

$ws.Range("A17").Value = "piauí"
$ws.Range("B17").Value = 2.71
$ws.Range("C17").Value = 12.04
$ws.Range("D17").Value = 4.83
$ws.Range("E17").Value = 10.35
$ws.Range("F17").Value = 9.27
$ws.Range("G17").Value = 12.73
$ws.Range("H17").Value = 15.36
$ws.Range("I17").Value = 19.53

$ws.Range("A18").Value = "ceará"
$ws.Range("B18").Value = 1.78
$ws.Range("C18").Value = 5.5
$ws.Range("D18").Value = 2.92
$ws.Range("E18").Value = 3.76
$ws.Range("F18").Value = 3.61
$ws.Range("G18").Value = 4.62
$ws.Range("H18").Value = 9.140000000000001
$ws.Range("I18").Value = 9.68

$ws.Range("A19").Value = "rio grande do norte"
$ws.Range("B19").Value = 6.87
$ws.Range("C19").Value = 11.53
$ws.Range("D19").Value = 5.87
$ws.Range("E19").Value = 12.4
$ws.Range("F19").Value = 11.47
$ws.Range("G19").Value = 10.95
$ws.Range("H19").Value = 18.82
$ws.Range("I19").Value = 13.89

$ws.Range("A20").Value = "paraíba"
$ws.Range("B20").Value = 4.74
$ws.Range("C20").Value = 13.42
$ws.Range("D20").Value = 5.83
$ws.Range("E20").Value = 4.49
$ws.Range("F20").Value = 8.140000000000001
$ws.Range("G20").Value = 7.31
$ws.Range("H20").Value = 17.31
$ws.Range("I20").Value = 18.25

$ws.Range("A21").Value = "pernambuco"
$ws.Range("B21").Value = 1.9
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 4.42
$ws.Range("F21").Value = 5.58
$ws.Range("G21").Value = 4.17
$ws.Range("H21").Value = 7
$ws.Range("I21").Value = 7.08

$ws.Range("A22").Value = "alagoas"
$ws.Range("B22").Value = 4.94
$ws.Range("C22").Value = 5.91
$ws.Range("D22").Value = 7.85
$ws.Range("E22").Value = 7.55
$ws.Range("F22").Value = 9.32
$ws.Range("G22").Value = 8.390000000000001
$ws.Range("H22").Value = 17.8
$ws.Range("I22").Value = 17.54

$ws.Range("A23").Value = "sergipe"
$ws.Range("B23").Value = 4.9
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 4.59
$ws.Range("E23").Value = 7.53
$ws.Range("F23").Value = 11.08
$ws.Range("G23").Value = 10.03
$ws.Range("H23").Value = 15.07
$ws.Range("I23").Value = 15.02

$ws.Range("A24").Value = "bahia"
$ws.Range("B24").Value = 1.88
$ws.Range("C24").Value = 5.43
$ws.Range("D24").Value = 2.58
$ws.Range("E24").Value = 4.03
$ws.Range("F24").Value = 3.06
$ws.Range("G24").Value = 3.4
$ws.Range("H24").Value = 8.19
$ws.Range("I24").Value = 7.47

$ws.Range("A25").Value = "sudeste"
$ws.Range("B25").Value = 0.66
$ws.Range("C25").Value = 2.42
$ws.Range("D25").Value = 1.14
$ws.Range("E25").Value = 1.65
$ws.Range("F25").Value = 1.69
$ws.Range("G25").Value = 1.21
$ws.Range("H25").Value = 2.78
$ws.Range("I25").Value = 2.89

$ws.Range("A26").Value = "minas gerais"
$ws.Range("B26").Value = 1.35
$ws.Range("C26").Value = 4.72
$ws.Range("D26").Value = 2.11
$ws.Range("E26").Value = 2.85
$ws.Range("F26").Value = 3.46
$ws.Range("G26").Value = 2.32
$ws.Range("H26").Value = 5.1
$ws.Range("I26").Value = 4.83

$ws.Range("A27").Value = "espírito santo"
$ws.Range("B27").Value = 2.68
$ws.Range("C27").Value = 7.2
$ws.Range("D27").Value = 3.57
$ws.Range("E27").Value = 5.36
$ws.Range("F27").Value = 5.38
$ws.Range("G27").Value = 6.62
$ws.Range("H27").Value = 12.55
$ws.Range("I27").Value = 9.73

$ws.Range("A28").Value = "rio de janeiro"
$ws.Range("B28").Value = 1.8
$ws.Range("C28").Value = 4.94
$ws.Range("D28").Value = 2.44
$ws.Range("E28").Value = 3.63
$ws.Range("F28").Value = 3.62
$ws.Range("G28").Value = 2.52
$ws.Range("H28").Value = 5.55
$ws.Range("I28").Value = 5.19

$ws.Range("A29").Value = "são paulo"
$ws.Range("B29").Value = 0.84
$ws.Range("C29").Value = 3.61
$ws.Range("D29").Value = 1.75
$ws.Range("E29").Value = 2.5
$ws.Range("F29").Value = 2.41
$ws.Range("G29").Value = 1.71
$ws.Range("H29").Value = 4.06
$ws.Range("I29").Value = 4.33

$ws.Range("A30").Value = "sul"
$ws.Range("B30").Value = 0.8100000000000001
$ws.Range("C30").Value = 4.92
$ws.Range("D30").Value = 1.56
$ws.Range("E30").Value = 2.13
$ws.Range("F30").Value = 2.04
$ws.Range("G30").Value = 1.87
$ws.Range("H30").Value = 3.94
$ws.Range("I30").Value = 4.15

$ws.Range("A31").Value = "paraná"
$ws.Range("B31").Value = 1.36
$ws.Range("C31").Value = 7.46
$ws.Range("D31").Value = 2.7
$ws.Range("E31").Value = 3.65
$ws.Range("F31").Value = 3.84
$ws.Range("G31").Value = 2.56
$ws.Range("H31").Value = 5.88
$ws.Range("I31").Value = 7.12

$ws.Range("A32").Value = "santa catarina"
$ws.Range("B32").Value = 1.85
$ws.Range("C32").Value = 9.960000000000001
$ws.Range("D32").Value = 4.53
$ws.Range("E32").Value = 4.76
$ws.Range("F32").Value = 3.69
$ws.Range("G32").Value = 4.73
$ws.Range("H32").Value = 9.58
$ws.Range("I32").Value = 9.16

$ws.Range("A33").Value = "rio grande do sul"
$ws.Range("B33").Value = 1.19
$ws.Range("C33").Value = 8.59
$ws.Range("D33").Value = 1.78
$ws.Range("E33").Value = 2.97
$ws.Range("F33").Value = 2.99
$ws.Range("G33").Value = 2.73
$ws.Range("H33").Value = 5.49
$ws.Range("I33").Value = 5.42

$ws.Range("A34").Value = "centro-oeste"
$ws.Range("B34").Value = 1.16
$ws.Range("C34").Value = 3.98
$ws.Range("D34").Value = 1.79
$ws.Range("E34").Value = 2.75
$ws.Range("F34").Value = 2.72
$ws.Range("G34").Value = 2.1
$ws.Range("H34").Value = 4.2
$ws.Range("I34").Value = 4.72

$ws.Range("A35").Value = "mato grosso do sul"
$ws.Range("B35").Value = 2.06
$ws.Range("C35").Value = 4.9
$ws.Range("D35").Value = 5
$ws.Range("E35").Value = 5.17
$ws.Range("F35").Value = 7.24
$ws.Range("G35").Value = 5.07
$ws.Range("H35").Value = 10.89
$ws.Range("I35").Value = 11.88

$ws.Range("A36").Value = "mato grosso"
$ws.Range("B36").Value = 2.7
$ws.Range("C36").Value = 10.83
$ws.Range("D36").Value = 3.64
$ws.Range("E36").Value = 7.66
$ws.Range("F36").Value = 6.16
$ws.Range("G36").Value = 4.31
$ws.Range("H36").Value = 10.15
$ws.Range("I36").Value = 12.56

$ws.Range("A37").Value = "goiás"
$ws.Range("B37").Value = 1.8
$ws.Range("C37").Value = 5.41
$ws.Range("D37").Value = 2.34
$ws.Range("E37").Value = 3.96
$ws.Range("F37").Value = 3.68
$ws.Range("G37").Value = 3.19
$ws.Range("H37").Value = 6.46
$ws.Range("I37").Value = 7.61

$ws.Range("A38").Value = "distrito federal"
$ws.Range("B38").Value = 3.15
$ws.Range("C38").Value = 8.789999999999999
$ws.Range("D38").Value = 4.71
$ws.Range("E38").Value = 6.3
$ws.Range("F38").Value = 6.44
$ws.Range("G38").Value = 5.05
$ws.Range("H38").Value = 7.95
$ws.Range("I38").Value = 8.220000000000001
